# edit.ps1 - apply LOT2028 worksheet restructuring (rows 13-25 -> rows 13-24)
# The original rows 13-25 get re-sequenced: column A (field labels) shifts up
# by one row relative to columns B/C (values) starting at row 13, a new
# "Semestral" value is introduced, and the final (now-duplicate) row 25 is removed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteValues = -4163
$xlPasteFormats = -4122

# --- Stage 1: snapshot every source cell (value + format) into a scratch area
# (rows 313-324) so later overwrites of rows 13-24 never clobber a value that
# is still needed as a source for another destination cell. ---
$stagePairs = @(
    @("A14", "A313"),
    @("A15", "A314"),
    @("B15", "B314"),
    @("C15", "C314"),
    @("A16", "A315"),
    @("B8", "B315"),
    @("C8", "C315"),
    @("A17", "A316"),
    @("B17", "B316"),
    @("C17", "C316"),
    @("A18", "A317"),
    @("A19", "A318"),
    @("B13", "B318"),
    @("C13", "C318"),
    @("A20", "A319"),
    @("B19", "B319"),
    @("C19", "C319"),
    @("A21", "A320"),
    @("B20", "B320"),
    @("C20", "C320"),
    @("A22", "A321"),
    @("B21", "B321"),
    @("C21", "C321"),
    @("A23", "A322"),
    @("B24", "B323"),
    @("C24", "C323"),
    @("B25", "B324"),
    @("C25", "C324")
)
foreach ($p in $stagePairs) {
    $src = $p[0]; $dst = $p[1]
    $ws.Range($src).Copy()
    $ws.Range($dst).PasteSpecial($xlPasteValues)
    $ws.Range($src).Copy()
    $ws.Range($dst).PasteSpecial($xlPasteFormats)
}

# --- Stage 2: the one brand-new string that has no prior occurrence in the sheet ---
$ws.Range("B13").Value2 = "Semestral"
$ws.Range("C13").Value2 = "Semestral"
$ws.Range("B19").Copy()
$ws.Range("B13").PasteSpecial($xlPasteFormats)
$ws.Range("C19").Copy()
$ws.Range("C13").PasteSpecial($xlPasteFormats)

# --- Stage 3: clear the B/C cells on rows that become single-column (label-only) rows ---
$ws.Range("B17,C17,B22,C22").Clear()

# --- Stage 4: copy the staged values+formats into their final rows 13-24 ---
$finalPairs = @(
    @("A313", "A13"),
    @("A314", "A14"),
    @("B314", "B14"),
    @("C314", "C14"),
    @("A315", "A15"),
    @("B315", "B15"),
    @("C315", "C15"),
    @("A316", "A16"),
    @("B316", "B16"),
    @("C316", "C16"),
    @("A317", "A17"),
    @("A318", "A18"),
    @("B318", "B18"),
    @("C318", "C18"),
    @("A319", "A19"),
    @("B319", "B19"),
    @("C319", "C19"),
    @("A320", "A20"),
    @("B320", "B20"),
    @("C320", "C20"),
    @("A321", "A21"),
    @("B321", "B21"),
    @("C321", "C21"),
    @("A322", "A22"),
    @("B323", "B23"),
    @("C323", "C23"),
    @("B324", "B24"),
    @("C324", "C24")
)
foreach ($p in $finalPairs) {
    $src = $p[0]; $dst = $p[1]
    $ws.Range($src).Copy()
    $ws.Range($dst).PasteSpecial($xlPasteValues)
    $ws.Range($src).Copy()
    $ws.Range($dst).PasteSpecial($xlPasteFormats)
}

# --- Stage 5: clear the scratch area entirely ---
$ws.Range("A313:C324").Clear()

# --- Stage 6: row heights for rows 13-24 ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30

# --- Stage 7: remove the now-obsolete row 25 (its content moved into row 24) ---
$ws.Rows.Item(25).Delete()

Write-Output "done"